# Applies the diff described in the commit to before.docx

$d = $word.ActiveDocument

# 1) "- графовые. Neo4j;" -> "- графовые: Neo4j;"
$d.Content.Find.Execute("- графовые. Neo4j;", $true, $false, $false, $false, $false,
                         $true, 1, $false, "- графовые: Neo4j;", 2)

# 2) Merge the "- " + "мультимодельные: OrientDB, ArangoDB. [" runs into one run,
#    and drop the period before "[": "- мультимодельные: OrientDB, ArangoDB. [" -> "- мультимодельные: OrientDB, ArangoDB ["
$d.Content.Find.Execute("- мультимодельные: OrientDB, ArangoDB. [", $true, $false, $false, $false, $false,
                         $true, 1, $false, "- мультимодельные: OrientDB, ArangoDB [", 2)

# 3) "]" (closing the Elmasri citation) -> "]." -- only touch the lone "]" run,
#    keep the distinctly-formatted "Elmasri" run untouched/separate.
$rng = $d.Content
$rng.Find.Execute("ArangoDB [Elmasri]")
$closeBracket = $d.Range($rng.End - 1, $rng.End)
$closeBracket.Text = "]."

# 4) Merge runs in the "Хранилища «ключ-значение»..." paragraph (pure concatenation, no text change)
$d.Content.Find.Execute("«ключа» и ассоциируемого с ним «значений», которое обычно является массивом данных. Подобные", $true, $false, $false, $false, $false,
                         $true, 1, $false, "«ключа» и ассоциируемого с ним «значений», которое обычно является массивом данных. Подобные", 2)

# 5) Merge runs in the "СУБД только с коммерческой лицензией..." paragraph (pure concatenation, no text change)
$d.Content.Find.Execute("поэтому СУБД Oracle в качестве СУБД для хранения знаний не рассматривается.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "поэтому СУБД Oracle в качестве СУБД для хранения знаний не рассматривается.", 2)

# 6) Merge "1." + "2" + " Современное состояние..." into one run
$d.Content.Find.Execute("1.2 Современное состояние области кластерного анализа данных", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1.2 Современное состояние области кластерного анализа данных", 2)

# 7) Merge runs in "Рассматривая область кластерного анализа данных..." paragraph (pure concatenation)
$d.Content.Find.Execute("кластеризации. Отдельные исследователи предлагают различные модели классификации, однако среди каждой из них можно выделить несколько основных направлений.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "кластеризации. Отдельные исследователи предлагают различные модели классификации, однако среди каждой из них можно выделить несколько основных направлений.", 2)

# 8) "Иерархические алгоритмы: аггломеративные, дивизимные (BIRCH, ...)" - drop trailing space before "("
$d.Content.Find.Execute("дивизимные (BIRCH, CURE, ROCK, Chameleon, Echidna)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "дивизимные (BIRCH, CURE, ROCK, Chameleon, Echidna)", 2)

# 9) "Разделяющие (K-means, ...)" - drop trailing space before "("
$d.Content.Find.Execute("Разделяющие (K-means, K-medoids, K-modes, PAM, CLARANS, CLARA,  FCM)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Разделяющие (K-means, K-medoids, K-modes, PAM, CLARANS, CLARA,  FCM)", 2)

# 10) "Плотностные (DBSCAN, ...)" - drop trailing space before "("
$d.Content.Find.Execute("Плотностные (DBSCAN, OPTICS, DBCLASD, DENCLUE)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Плотностные (DBSCAN, OPTICS, DBCLASD, DENCLUE)", 2)

# 11) "Сеточные (Wave-Cluster, ...)" - drop trailing space before "("
$d.Content.Find.Execute("Сеточные (Wave-Cluster, STING, CLIQUE, OptiGrid)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Сеточные (Wave-Cluster, STING, CLIQUE, OptiGrid)", 2)

# 12) "Моделируемые (EM, COBWEB, CLASSIT, SOMs)" - drop trailing space before "("
$d.Content.Find.Execute("Моделируемые (EM, COBWEB, CLASSIT, SOMs)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Моделируемые (EM, COBWEB, CLASSIT, SOMs)", 2)
